$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header row text changes (A1:H1), reordering + rewording columns
# ---------------------------------------------------------------
# A1: 單位名稱 -> 所屬一級單位 (rich text: "所屬一級" red, "單位" black)
$ws.Range("A1").Value = "所屬一級單位"
# B1: 系所部門 -> 所屬系所部門 (rich text: "所屬" red, "系所部門" black)
$ws.Range("B1").Value = "所屬系所部門"
# C1: 姓名 (unchanged text)
$ws.Range("C1").Value = "姓名"
# D1: 身分 -> 身分 (學士、碩士或博士班）
$ws.Range("D1").Value = "身分 (學士、碩士或博士班）"
# E1: 國籍 moved up from the end of the row to column E
$ws.Range("E1").Value = "國籍"
# F1: 開始時間 (unchanged text, shifted from column E)
$ws.Range("F1").Value = "開始時間"
# G1: 結束時間 (unchanged text, shifted from column F)
$ws.Range("G1").Value = "結束時間"
# H1: 備註 (unchanged text, shifted from column G)
$ws.Range("H1").Value = "備註"

# ---------------------------------------------------------------
# Header formatting: bold, vertically centered for the whole row,
# with D1 (身分...) recoloured blue instead of the default black.
# ---------------------------------------------------------------
$ws.Range("A1:H1").Font.Bold = $true
$ws.Range("A1:H1").VerticalAlignment = -4108
$ws.Range("D1").Font.Color = 16711680

# ---------------------------------------------------------------
# Two-colour rich text runs for A1 and B1: a red-highlighted
# leading phrase ("所屬一級" / "所屬") followed by a black phrase.
# ---------------------------------------------------------------
$ws.Range("A1").Characters(1, 4).Font.Color = 255
$ws.Range("A1").Characters(1, 4).Font.Bold = $true
$ws.Range("A1").Characters(5, 2).Font.Color = 0
$ws.Range("A1").Characters(5, 2).Font.Bold = $true

$ws.Range("B1").Characters(1, 2).Font.Color = 255
$ws.Range("B1").Characters(1, 2).Font.Bold = $true
$ws.Range("B1").Characters(3, 4).Font.Color = 0
$ws.Range("B1").Characters(3, 4).Font.Bold = $true

# ---------------------------------------------------------------
# Column widths for A, B and D
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 97 / 7
$ws.Columns("B").ColumnWidth = 95 / 7
$ws.Columns("D").ColumnWidth = 188 / 7
